$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'track tights'
$ws.Cells.Item(2, 1).Value = 'wintergear compression men'
$ws.Cells.Item(3, 1).Value = 'basketball training tights'
$ws.Cells.Item(4, 1).Value = 'basketball protector'
$ws.Cells.Item(5, 1).Value = 'winter leggings for men'
$ws.Cells.Item(6, 1).Value = 'kids tights with knee pads'
$ws.Cells.Item(7, 1).Value = 'ropa de monta?a hombre'
$ws.Cells.Item(8, 1).Value = 'men workout leggings nike'
$ws.Cells.Item(9, 1).Value = 'boys white knee pads basketball'
$ws.Cells.Item(10, 1).Value = 'coyote brown pants with knee pads'
$ws.Cells.Item(11, 1).Value = 'black basketball knee pads nike'
$ws.Cells.Item(12, 1).Value = 'basketball knee pads womens'
$ws.Cells.Item(13, 1).Value = 'cold weather workout pants for men'
$ws.Cells.Item(14, 1).Value = 'mens under armour long underwear pants'
$ws.Cells.Item(15, 1).Value = 'calf compression leggings men'
$ws.Cells.Item(16, 1).Value = 'tights for men nike'
$ws.Cells.Item(17, 1).Value = 'swimming pants for men'
$ws.Cells.Item(18, 1).Value = 'copper compression pants for men'
$ws.Cells.Item(19, 1).Value = 'tesla thermal pants'
$ws.Cells.Item(20, 1).Value = 'hex knee sleeve'
$ws.Cells.Item(21, 1).Value = 'thermal nike'
$ws.Cells.Item(22, 1).Value = 'running base layer men'
$ws.Cells.Item(23, 1).Value = 'long spandex men'
$ws.Cells.Item(24, 1).Value = 'cold gear compression pants men'
$ws.Cells.Item(25, 1).Value = 'soccer winter gear'
$ws.Cells.Item(26, 1).Value = 'kids soccor gear'
$ws.Cells.Item(27, 1).Value = 'mens pants with knee pads'
$ws.Cells.Item(28, 1).Value = 'mens thermal baselayer'
$ws.Cells.Item(29, 1).Value = 'knee pad pants men'
$ws.Cells.Item(30, 1).Value = 'adidas compression pants men'
$ws.Cells.Item(31, 1).Value = 'adult clothing protector'
$ws.Cells.Item(32, 1).Value = 'adult football knee pads'
$ws.Cells.Item(33, 1).Value = 'adult football pants with pads'
$ws.Cells.Item(34, 1).Value = 'athletic knee pads'
$ws.Cells.Item(35, 1).Value = 'athletic tights men'
$ws.Cells.Item(36, 1).Value = 'baseball knee pad'
$ws.Cells.Item(37, 1).Value = 'baseball knee pads'
$ws.Cells.Item(38, 1).Value = 'baseball pants adult small'
$ws.Cells.Item(39, 1).Value = 'baskerball tights'
$ws.Cells.Item(40, 1).Value = 'basketball compression knee pads'
$ws.Cells.Item(41, 1).Value = 'basketball compression pants youth with knee pads'
$ws.Cells.Item(42, 1).Value = 'basketball for youth'
$ws.Cells.Item(43, 1).Value = 'basketball hip pads'
$ws.Cells.Item(44, 1).Value = 'basketball leggings men'
$ws.Cells.Item(45, 1).Value = 'basketball pants for women'
$ws.Cells.Item(46, 1).Value = 'basketball pants men'
$ws.Cells.Item(47, 1).Value = 'basketball tights boys youth'
$ws.Cells.Item(48, 1).Value = 'basketball tights for men mcdavid'
$ws.Cells.Item(49, 1).Value = 'basketball tights with pads for boys'
$ws.Cells.Item(50, 1).Value = 'best basketball knee pads'
$ws.Cells.Item(51, 1).Value = 'big knee pads'
$ws.Cells.Item(52, 1).Value = 'big man knee pads'
$ws.Cells.Item(53, 1).Value = 'bjj leggings'
$ws.Cells.Item(54, 1).Value = 'black knee pads'
$ws.Cells.Item(55, 1).Value = 'black knee pads for volleyball'
$ws.Cells.Item(56, 1).Value = 'black leggings design'
$ws.Cells.Item(57, 1).Value = 'black mesh leggings capri'
$ws.Cells.Item(58, 1).Value = 'black youth baseball pants'
$ws.Cells.Item(59, 1).Value = 'boys athletic tights basketball'
$ws.Cells.Item(60, 1).Value = 'boys basketball knee pads mcdavid'
$ws.Cells.Item(61, 1).Value = 'boys compression leggings'
$ws.Cells.Item(62, 1).Value = 'boys compression leggings youth'
$ws.Cells.Item(63, 1).Value = 'boys compression tights'
$ws.Cells.Item(64, 1).Value = 'break away basketball pants'
$ws.Cells.Item(65, 1).Value = 'capri leggings medium'
$ws.Cells.Item(66, 1).Value = 'capri mens'
$ws.Cells.Item(67, 1).Value = 'capri tights for men'
$ws.Cells.Item(68, 1).Value = 'capris leggings'
$ws.Cells.Item(69, 1).Value = 'cold gear for football'
$ws.Cells.Item(70, 1).Value = 'cold weather panta'
$ws.Cells.Item(71, 1).Value = 'compression for knee'
$ws.Cells.Item(72, 1).Value = 'compression gear'
$ws.Cells.Item(73, 1).Value = 'compression men pants'
$ws.Cells.Item(74, 1).Value = 'compression pants big and tall men'
$ws.Cells.Item(75, 1).Value = 'compression pants padded knees basketball'
$ws.Cells.Item(76, 1).Value = 'compression shorts 3 4 length men'
$ws.Cells.Item(77, 1).Value = 'compression tights for men'
$ws.Cells.Item(78, 1).Value = 'cycling pants'
$ws.Cells.Item(79, 1).Value = 'dark purple basketball knee pads'
$ws.Cells.Item(80, 1).Value = 'elbow knee pads youth'
$ws.Cells.Item(81, 1).Value = 'excersize equipment for men'
$ws.Cells.Item(82, 1).Value = 'football 3 4 tights'
$ws.Cells.Item(83, 1).Value = 'football knee pads'
$ws.Cells.Item(84, 1).Value = 'football leg pads'
$ws.Cells.Item(85, 1).Value = 'football pants'
$ws.Cells.Item(86, 1).Value = 'football pants adult black'
$ws.Cells.Item(87, 1).Value = 'football pants youth'
$ws.Cells.Item(88, 1).Value = 'g form knee pads youth'
$ws.Cells.Item(89, 1).Value = 'gel knee pads'
$ws.Cells.Item(90, 1).Value = 'gel knee pads for men'
$ws.Cells.Item(91, 1).Value = 'girl knee pads'
$ws.Cells.Item(92, 1).Value = 'girls basketball knee pads'
$ws.Cells.Item(93, 1).Value = 'girls tights with knee pads'
$ws.Cells.Item(94, 1).Value = 'girls volleyball knee pads'
$ws.Cells.Item(95, 1).Value = 'girls volleyball knee pads youth'
$ws.Cells.Item(96, 1).Value = 'girls youth volleyball knee pads'
$ws.Cells.Item(97, 1).Value = 'gym pants for men'
$ws.Cells.Item(98, 1).Value = 'happy knees'
$ws.Cells.Item(99, 1).Value = 'hex foam knee pads'
$ws.Cells.Item(100, 1).Value = 'hex knee pads for basketball'
